$wb = $excel.ActiveWorkbook

# The previously-active sheet ("ciudad_distinta") keeps its own cursor
# position, just no longer on the selected/active tab - update its
# remembered selection before we move focus away from it.
$prevSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$prevSheet.Range("E14").Select()

# Add the new "nombre" worksheet after the last existing sheet (ciudad_distinta).
$newSheet = $wb.Worksheets.Add($null, $prevSheet)
$newSheet.Name = "nombre"

# Populate the new sheet with the "nombre" list.
$newSheet.Range("A1").Value = "nombre"
$newSheet.Range("A2").Value = "Ana"
$newSheet.Range("A3").Value = "Lucía"
$newSheet.Range("A4").Value = "María"

# The last two entries wrap text, matching the new style added to cellXfs.
$newSheet.Range("A3:A4").WrapText = $true

# Make the new sheet the active / selected sheet and set its selection.
$newSheet.Activate()
$newSheet.Range("B12").Select()
